$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column D (the second attendance-check column, rows 1-17) into the
# new column E so the attendance table gains a matching scrollable column.
# Using Range.Copy (instead of re-typing the values) preserves the original
# cell types/formatting exactly - e.g. the "2025-06-18" header and the
# check-mark text stay literal text instead of being re-interpreted as dates.
$src = $ws.Range("D1:D17")
$dst = $ws.Range("E1:E17")
$src.Copy($dst)
